$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for the 2022-Q3 summary, pushing
#    the existing quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 32
$total.Cells.Item(2,4).Value = 8.640000000000001

# Clear the formatting the row-insert borrowed from row 3, then restore
# the index-column style (matches the rest of column A) by copying it
# from the (now shifted) row below.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" worksheet right after "总计" holding
#    the quarter's fund-holder detail table.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Header row
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$hdrRange = $q3.Range("B1:H1")
$hdrRange.Font.Bold = $true
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4160
$hdrRange.Borders.LineStyle = 1

# Data rows: index(A,n) code(B,text) name(C,text) size(D,text) position(E,text)
#            ratio(F,text) value(G,text) rank(H,n)
$rows = @(
    @(0,  "001955", "中欧养老产业混合A",              "22.51", "92.70", "7.91", "1.7805", 8),
    @(1,  "008120", "万家自主创新混合A",              "22.86", "93.71", "7.44", "1.7008", 4),
    @(2,  "010429", "中欧睿见混合A",                  "18.74", "91.22", "6.34", "1.1881", 10),
    @(3,  "005311", "万家经济新动能混合A",            "10.59", "93.55", "7.67", "0.8123", 5),
    @(4,  "005421", "中欧嘉泽灵活配置混合",           "6.99",  "83.96", "6.85", "0.4788", 2),
    @(5,  "005312", "万家经济新动能混合C",            "5.65",  "93.55", "7.67", "0.4334", 5),
    @(6,  "011708", "中欧嘉益一年混合A",              "4.44",  "93.35", "7.66", "0.3401", 2),
    @(7,  "012778", "中欧养老产业混合C",              "2.80",  "92.70", "7.91", "0.2215", 8),
    @(8,  "008633", "万家科技创新混合A",              "2.44",  "93.37", "8.25", "0.2013", 2),
    @(9,  "008121", "万家自主创新混合C",              "2.24",  "93.71", "7.44", "0.1667", 4),
    @(10, "004616", "中欧电子信息产业沪港深股票A",    "5.01",  "92.97", "3.32", "0.1663", 7),
    @(11, "005416", "鹏华尊惠18个月定期开放混合A",    "7.53",  "37.52", "1.92", "0.1446", 4),
    @(12, "011709", "中欧嘉益一年混合C",              "1.87",  "93.35", "7.66", "0.1432", 2),
    @(13, "008634", "万家科技创新混合C",              "1.71",  "93.37", "8.25", "0.1411", 2),
    @(14, "550001", "信诚四季红混合",                 "4.79",  "83.57", "2.84", "0.1360", 10),
    @(15, "005763", "中欧电子信息产业沪港深股票C",    "3.88",  "92.97", "3.32", "0.1288", 7),
    @(16, "673141", "西部利得景程灵活配置混合A",      "1.36",  "85.50", "5.63", "0.0766", 4),
    @(17, "007854", "光大保德信景气先锋混合",         "1.50",  "73.83", "4.95", "0.0742", 8),
    @(18, "009667", "鹏华安庆混合A",                  "2.12",  "38.12", "1.95", "0.0413", 4),
    @(19, "009230", "鹏华安和混合A",                  "1.80",  "39.30", "1.91", "0.0344", 4),
    @(20, "011284", "中信保诚龙腾精选混合",           "1.12",  "84.65", "2.85", "0.0319", 10),
    @(21, "673143", "西部利得景程灵活配置混合C",      "0.54",  "85.50", "5.63", "0.0304", 4),
    @(22, "011572", "鹏华安荣混合A",                  "1.50",  "39.92", "1.98", "0.0297", 4),
    @(23, "003165", "鹏华弘嘉灵活配置混合A",          "0.82",  "91.09", "3.53", "0.0289", 5),
    @(24, "009668", "鹏华安庆混合C",                  "1.32",  "38.12", "1.95", "0.0257", 4),
    @(25, "001723", "华商新动力灵活配置混合",         "0.70",  "84.40", "3.25", "0.0228", 9),
    @(26, "006209", "中信保诚新蓝筹灵活配置混合",     "0.52",  "80.11", "2.71", "0.0141", 10),
    @(27, "009231", "鹏华安和混合C",                  "0.73",  "39.30", "1.91", "0.0139", 4),
    @(28, "005417", "鹏华尊惠18个月定期开放混合C",    "0.58",  "37.52", "1.92", "0.0111", 4),
    @(29, "003166", "鹏华弘嘉灵活配置混合C",          "0.25",  "91.09", "3.53", "0.0088", 5),
    @(30, "011573", "鹏华安荣混合C",                  "0.34",  "39.92", "1.98", "0.0067", 4),
    @(31, "015481", "中欧睿见混合C",                  "0.08",  "91.22", "6.34", "0.0051", 10)
)

$textCols = 2,3,4,5,6,7   # B..G hold text values

foreach ($r in $rows) {
    $rowNum = [int]$r[0] + 2
    $q3.Cells.Item($rowNum, 1).Value = $r[0]
    foreach ($c in $textCols) {
        $cell = $q3.Cells.Item($rowNum, $c)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$r[$c - 1]
    }
    $q3.Cells.Item($rowNum, 8).Value = $r[7]
}

# Match the index column's look (bold / centered / bordered) to the rest
# of the workbook's "index" columns.
$idxRange = $q3.Range("A2:A33")
$idxRange.Font.Bold = $true
$idxRange.HorizontalAlignment = -4108
$idxRange.VerticalAlignment = -4160
$idxRange.Borders.LineStyle = 1

$q3.Application.CutCopyMode = $false
